$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-converted to a number by Excel
# must be forced to Text format first so they remain stored as text (matching the
# original inline-string / shared-string cell type).
$textCells = @("D5","D6","D8","D9","D10","D11","D12","D15","D17","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D38","D39","D40","D41","D42","D43","D45","D47","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '26.224.89'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '1.656.87'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  -0.67%  '
$ws.Range('D5').Value = '218.97'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').Value = '0.5239'
$ws.Range('E6').Value = '  -2.02%  '
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('D8').Value = '0.2642'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').Value = '0.06310'
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('D10').Value = '20.62'
$ws.Range('E10').Value = '  -1.40%  '
$ws.Range('D11').Value = '0.07777'
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').Value = '4.513'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').Value = '1.644.46'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').Value = '1.886.89'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').Value = '0.5629'
$ws.Range('E15').Value = '  +1.53%  '
$ws.Range('D16').Value = '0.0₅8058'
$ws.Range('E16').Value = '  -1.46%  '
$ws.Range('D17').Value = '65.22'
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').Value = '26.231.22'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('D20').Value = '4.717'
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').Value = '193.57'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('D22').Value = '10.24'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = '6.031'
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').Value = '145.16'
$ws.Range('D26').Value = '0.1207'
$ws.Range('E26').Value = '  -1.46%  '
$ws.Range('D27').Value = '7.227'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').Value = '16.01'
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('D29').Value = '1.497'
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('D30').Value = '0.05621'
$ws.Range('E30').Value = '  -4.02%  '
$ws.Range('D31').Value = '1.279'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('D32').Value = '3.486'
$ws.Range('E32').Value = '  -2.62%  '
$ws.Range('D33').Value = '3.369'
$ws.Range('E33').Value = '  +2.49%  '
$ws.Range('D34').Value = '1.600'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('D35').Value = '2.803'
$ws.Range('E35').Value = '  -1.03%  '
$ws.Range('D36').Value = '0.9445'
$ws.Range('E36').Value = '  -2.68%  '
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('D38').Value = '0.5746'
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('D39').Value = '0.01603'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').Value = '5.975'
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('D41').Value = '2.568'
$ws.Range('E41').Value = '  -0.68%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.8470'
$ws.Range('E42').Value = '  -2.02%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '1.004'
$ws.Range('E43').Value = '  -0.72%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.040.65'
$ws.Range('E44').Value = '  -3.52%  '
$ws.Range('D45').Value = '102.11'
$ws.Range('E45').Value = '  -2.11%  '
$ws.Range('D46').Value = '1.797.94'
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('D47').Value = '58.30'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('E48').Value = '  +1.33%  '
$ws.Range('E49').Value = '  -1.28%  '
$ws.Range('D50').Value = '0.05313'
$ws.Range('E50').Value = '  +2.80%  '
$ws.Range('D51').Value = '8.059'
$ws.Range('E51').Value = '  -0.01%  '
